$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 13
$ws1.Range("F6").Value = 3175
$ws1.Range("F7").Value = 2085
$ws1.Range("F9").Value = 149
$ws1.Range("F12").Value = 1059
$ws1.Range("F13").Value = 86

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 13
$ws4.Range("F6").Value = 3175
$ws4.Range("F7").Value = 2085
$ws4.Range("F10").Value = 149
$ws4.Range("F13").Value = 1059
$ws4.Range("F14").Value = 86
